$wb = $excel.ActiveWorkbook

# --- Sheet "CompanyList" (sheet1.xml): reorder company rows and append "ABC" ---
$ws1 = $wb.Worksheets.Item("CompanyList")
$ws1.Range("A3").Value = "Just Dial"
$ws1.Range("A4").Value = "DLF"
$ws1.Range("A5").Value = "Aditya Birla F"
$ws1.Range("A6").Value = "Ceat"
$ws1.Range("A7").Value = "ABC"

# --- Sheet "Sheet1" (sheet2.xml): add column B mirroring column A, shifted by one row ---
$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Range("B1").Value = "Aditya Birla F"
$ws2.Range("B2").Value = "Ceat"
$ws2.Range("B3").Value = "Just Dial"
$ws2.Range("B4").Value = "DFL"

# Column widths for the newly visible columns on "Sheet1"
$ws2.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws2.Columns.Item(2).ColumnWidth = 14.833333333333334

# Selections: set "Sheet1" selection first, then "CompanyList" so CompanyList stays the active tab
$ws2.Range("B1").Select()
$ws1.Range("C9").Select()
